$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal TEXT (avoid Excel auto-converting
# numeric-looking strings like "1.00" or "578.20" into numbers), while
# keeping the cell on the default/unstyled format (same as source file).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "65.092.01"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.478.12"
$ws.Range("E3").Value = "  +3.40%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
Set-TextValue $ws.Range("D5") "578.20"
$ws.Range("E5").Value = "  +3.66%  "

# Row 6
Set-TextValue $ws.Range("D6") "160.32"
$ws.Range("E6").Value = "  +4.67%  "

# Row 7
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue $ws.Range("D8") "0.607"
$ws.Range("E8").Value = "  +13.01%  "

# Row 9
Set-TextValue $ws.Range("D9") "3.474.02"
$ws.Range("E9").Value = "  +3.36%  "

# Row 10
Set-TextValue $ws.Range("D10") "7.25"
$ws.Range("E10").Value = "  -1.22%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.125"
$ws.Range("E11").Value = "  +4.16%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.444"
$ws.Range("E12").Value = "  +3.67%  "

# Row 13
Set-TextValue $ws.Range("D13") "4.070.92"
$ws.Range("E13").Value = "  +3.10%  "

# Row 14
$ws.Range("E14").Value = "  +1.21%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.0000192"
$ws.Range("E15").Value = "  +3.25%  "

# Row 16
Set-TextValue $ws.Range("D16") "28.50"
$ws.Range("E16").Value = "  +7.22%  "

# Row 17
Set-TextValue $ws.Range("D17") "65.058.24"
$ws.Range("E17").Value = "  +3.52%  "

# Row 18
Set-TextValue $ws.Range("D18") "3.481.13"
$ws.Range("E18").Value = "  +3.09%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.43"
$ws.Range("E19").Value = "  +4.26%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.28"
$ws.Range("E20").Value = "  +2.93%  "

# Row 21
Set-TextValue $ws.Range("D21") "381.45"
$ws.Range("E21").Value = "  +2.12%  "

# Row 22
Set-TextValue $ws.Range("D22") "8.20"
$ws.Range("E22").Value = "  +3.72%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.549"
$ws.Range("E23").Value = "  +5.21%  "

# Row 24
Set-TextValue $ws.Range("D24") "72.91"
$ws.Range("E24").Value = "  +2.68%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.00"
$ws.Range("E25").Value = "  +0.59%  "

# Row 26
Set-TextValue $ws.Range("D26") "0.0000119"
$ws.Range("E26").Value = "  +5.42%  "

# Row 27
Set-TextValue $ws.Range("D27") "9.95"
$ws.Range("E27").Value = "  +6.68%  "

# Row 28
$ws.Range("E28").Value = "  +2.42%  "

# Row 29
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
Set-TextValue $ws.Range("D30") "1.49"
$ws.Range("E30").Value = "  +12.59%  "

# Row 31
Set-TextValue $ws.Range("D31") "6.18"
$ws.Range("E31").Value = "  +3.51%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.05"
$ws.Range("E32").Value = "  +4.48%  "

# Row 33
Set-TextValue $ws.Range("D33") "23.57"
$ws.Range("E33").Value = "  +2.99%  "

# Row 34
Set-TextValue $ws.Range("D34") "7.23"
$ws.Range("E34").Value = "  +8.79%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.57"
$ws.Range("E35").Value = "  +9.83%  "

# Row 36
Set-TextValue $ws.Range("D36") "161.23"
$ws.Range("E36").Value = "  +1.83%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.91"

# Row 38
Set-TextValue $ws.Range("D38") "3.006.73"
$ws.Range("E38").Value = "  +4.26%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.0772"
$ws.Range("E39").Value = "  +1.90%  "

# Row 40
Set-TextValue $ws.Range("D40") "26.90"
$ws.Range("E40").Value = "  +0.86%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0322"
$ws.Range("E41").Value = "  +2.90%  "

# Row 42
Set-TextValue $ws.Range("D42") "4.54"
$ws.Range("E42").Value = "  +6.90%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D43") "42.53"
$ws.Range("E43").Value = "  +4.44%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D44") "6.44"
$ws.Range("E44").Value = "  +0.20%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.777"
$ws.Range("E45").Value = "  +4.99%  "

# Row 46
Set-TextValue $ws.Range("D46") "25.55"
$ws.Range("E46").Value = "  +12.05%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.10"
$ws.Range("E47").Value = "  +4.89%  "

# Row 48
Set-TextValue $ws.Range("D48") "320.35"
$ws.Range("E48").Value = "  +13.37%  "

# Row 49
Set-TextValue $ws.Range("D49") "6.73"
$ws.Range("E49").Value = "  +7.10%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.109"
$ws.Range("E50").Value = "  +7.62%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.19"
$ws.Range("E51").Value = "  +5.03%  "
